$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: BD01 arm-root link gets repointed and its offset id bumped ---
$ws.Range("C9").Value = 200
$ws.Range("D9").Value = "BD01"

# --- Left Arm IK block (rows 43-48): mirror the existing Right Arm block
#     (rows 33-38) so every border/fill carries over exactly, then patch
#     the handful of cells that actually differ for the left side. ---
$ws.Range("A33:AE38").Copy($ws.Range("A43:AE48"))

$leftArmRows = @(
    @{ row = 43; A = 1400; B = "LAEE01"; C = 402; D = "LA03";  E = -2;    L = -1; AE = "Left Wrist Rot Obj" },
    @{ row = 44; A = 2400; B = "LATA01";                       E = -3;    L = -1; AE = "Left Wrist Tar" },
    @{ row = 45; A = 1401; B = "LAEE02"; C = 401; D = "LA02";  E = -2;    L = -1; AE = "Left elbow Obj" },
    @{ row = 46; A = 2401; B = "LATA02";                       E = -2.7;  L = -1; AE = "Left elbow Tar" },
    @{ row = 47; A = 1402; B = "LAEE03"; C = 402; D = "LA03";  E = -2;    L = -1; AE = "Left Wrist Pos Obj" },
    @{ row = 48; A = 2402; B = "LATA03";                       E = -2.2;  L = -1; AE = "Left Wrist Pos Tar" }
)

foreach ($r in $leftArmRows) {
    $row = $r.row
    $ws.Cells.Item($row, 1).Value = $r.A          # A - id
    $ws.Cells.Item($row, 2).Value = $r.B          # B - name
    if ($r.ContainsKey("C")) { $ws.Cells.Item($row, 3).Value = $r.C }   # C - parent id
    if ($r.ContainsKey("D")) { $ws.Cells.Item($row, 4).Value = $r.D }   # D - parent name
    $ws.Cells.Item($row, 5).Value = $r.E          # E - offset
    $ws.Cells.Item($row, 12).Value = $r.L         # L - rot axis x (mirrored sign)
    $ws.Cells.Item($row, 31).Value = $r.AE        # AE - comment
}

# row 48's offset (-2.2) must keep the same float precision quirk ("-2.2000000000000002")
# that mirrors row 38's (2.2000000000000002) floating point representation
$ws.Cells.Item(48, 5).Value = -2.2000000000000002

# --- View state: zoom out and move the cursor before saving, as in the source edit ---
$win = $excel.ActiveWindow
$win.Zoom = 55
$ws.Range("D9").Select()
